$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("log_data")

# New data for the past 3 days
$ws.Range("A67").Value = 43974
$ws.Range("B67").Value = 2627

$ws.Range("A68").Value = 43975
$ws.Range("B68").Value = 2648

$ws.Range("A69").Value = 43976
$ws.Range("B69").Value = 2676

# Fill the formula down the range so it is stored as one shared formula group
$ws.Range("C67:C69").Formula = "=LOG10(B67)"

# Match date formatting style used in column A (style index 1 -> yyyy-mm-dd)
$ws.Range("A67:A69").NumberFormat = "yyyy\-mm\-dd;@"

# Update the selection to match the new focus (F68)
$ws.Range("F68").Select()
